$d = $word.ActiveDocument

# The second paragraph holds a Word field (fldChar begin / instrText* / fldChar end)
# encoding an M2Doc expression: " m:'Mona_Lisa.jpg'.asImage().setWidth(100) ".
# The parser was switched to TokenIteratorFieldRewriterSplit, which rewrites the
# field back into plain literal text runs "{" ... "}" instead of a Word field,
# while preserving the run-level coloring (w:rPr/w:color) on the expression body.

$fld = $d.Fields(1)
$target = $fld.Code.Paragraphs(1)
$targetRange = $target.Range

# Single-quoted here-string: no interpolation needed, so the apostrophes in
# the M2Doc expression ('Mona_Lisa.jpg') can be written literally.
$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>'</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>Mona_Lisa</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.jpg</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>'.asImage()</w:t></w:r><w:r><w:rPr><w:color w:val="E36C0A" w:themeColor="accent6" w:themeShade="BF"/></w:rPr><w:t>.setWidth(100)</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$targetRange.InsertXML($xml)
